# Plantilla "asignaturasPlantilla.xlsx" update:
#  - Remove the merged "ASIGNATURAS" title row.
#  - Promote the column-header row (Codigo / Acronimo / Nombre / Titulacion)
#    up to row 1, and append a new "Curso" header in column E.
#  - Keep the leftover formatted marker cell (previously D11) which now
#    lands on D10 after the row shift, and add a matching one at E4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the merge on the old title row, then delete that whole row so
# everything below (the real header + the stray formatted cell) shifts up.
$ws.Range("A1:D1").UnMerge()
$ws.Rows("1").Delete()

# New "Curso" header in column E, matching the look of the other headers.
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("E1").Value = "Curso"

# Mirror the lone formatted cell (now D10) onto E4, matching the new layout.
$ws.Range("D10").Copy()
$ws.Range("E4").PasteSpecial(-4122)

$ws.Range("E4").Select() | Out-Null
